$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row numbers (A column) for the new module-7 entries
$ws.Range("A47").Value = 7
$ws.Range("A48").Value = 7
$ws.Range("A49").Value = 7
$ws.Range("A50").Value = 7
$ws.Range("A51").Value = 7

# Populate cells in the exact order the strings were authored so the
# shared-string table is built in the same sequence as the target file.
$ws.Range("C47").Value = "Inference in a frequentist framework draws conclusions from sample data by conceiving of this specific 'experiment' or sample as only one of thousands of possible experiments/samples, each capable of producing statistically independent results. Thus our inference is based on the probability of a given parameter (e.g. from one sample or experiment) arising in relation to all other (random) possibilities."
$ws.Range("B48").Value = "Bayesian Inference"
$ws.Range("B47").Value = "Frequentist Inference"
$ws.Range("C48").Value = "Bayesian is a process of using observed data to update prior beliefs. Typically parameters are assumed to be random variables arising from a distribution (e.g. rather than a discrete and solitary truth)."
$ws.Range("B49").Value = "Prior"
$ws.Range("C49").Value = "In Bayesian inference, the 'prior' is a formalized statement of the probability of a parameter, as stated before we see the data."
$ws.Range("B50").Value = "Posterior"
$ws.Range("C50").Value = "In Bayesian inference, the 'posterior' is a formalized statement about the updated belief of the value of a parameter, conditional on the data (the likelihood) and the prior."
$ws.Range("B51").Value = "Conditional auto-regressive (CAR) "
$ws.Range("C51").Value = "The CAR is a common prior for spatial disease mapping, particularly in a Bayesian framework. A CAR prior suggests that the value for a given area can be estimated CONDITIONAL ON the level of neighboring values. "

# Update the view so row 47 onward is visible and the active cell sits on the new blank row.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("A52").Select()
